$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Prefix with an apostrophe so Excel stores these as text (shared strings)
# rather than converting them to numeric values, matching the source data.
$ws.Range("B10").Value = "'73.62"
$ws.Range("C10").Value = "'23.74"
$ws.Range("D10").Value = "'97.36"
